# "Add files via upload" - restructure PeopleList sheet:
#   A=uid, B=original name (raw), C=name (extracted), D=introducer (extracted)
# Insert a new column before the old "name" column (B), which pushes the
# old B ("name") to C and old C ("introducer") to D. Then set the new B1
# header, and replace C/D's per-row placeholder cells with extraction
# formulas driven off column B (only needed on row 2, matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B - shifts "name" -> C, "introducer" -> D
$ws.Columns.Item(2).Insert()

# New header text for the inserted column
$ws.Range("B1").Value = "original name"

# The inserted column is style-less below the header; give it back the
# same unlocked text-entry style the old "name" column used to carry.
$ws.Range("B2:B28").Style = $ws.Range("C2").Style

# Clear out the old per-row empty placeholder cells that used to live in
# column C (now shifted to C/D) below row 2 - the new layout only keeps
# data entry in column B, with C/D driven by formulas on row 2.
$ws.Range("C3:D28").ClearContents()
$ws.Range("C3:D28").ClearFormats()
$ws.Range("C29:D30").Delete()

# Formulas that derive "name" and "introducer" out of the raw "original name"
$ws.Range("C2").Formula = '=LEFT(B2, FIND("(", B2) - 1)'
$ws.Range("D2").Formula = '=IFERROR(MID(B2, FIND("(", B2), LEN(B2) - FIND("(", B2) + 1), "")'

# Column widths
$ws.Columns.Item(1).ColumnWidth = 22.140625
$ws.Columns.Item(2).ColumnWidth = 36.28515625
$ws.Columns.Item(3).ColumnWidth = 16.5703125
$ws.Columns.Item(4).ColumnWidth = 23.140625

# Row heights / view tweaks
$ws.Rows.Item(1).RowHeight = 27.75
$ws.Range("D2").Select()

$wb.Save()
